$wb = $excel.ActiveWorkbook

# This script applies a scheduled market-price data refresh to the per-job
# "Leve Profit" tables (columns H:N) across all 8 job sheets, mirroring an
# external price-fetch run. Cell values are written directly; a couple of
# rows additionally gain/lose their Profit(NQ)/Profit(HQ) cell depending on
# whether that HQ/NQ branch is priced this run.

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 897.4706
$ws.Range("J12").Value = 3549.75
$ws.Range("L12").Value = 3549.75
$ws.Range("N12").Value = -3889.75
# Row 45
$ws.Range("H45").Value = 6470.6
$ws.Range("I45").Value = 700
$ws.Range("J45").Value = 7913.25
$ws.Range("K45").Value = 2100
$ws.Range("L45").Value = 23739.75
$ws.Range("M45").Value = -1908
$ws.Range("N45").Value = -24123.75
# Row 137
$ws.Range("H137").Value = 2591.56
$ws.Range("I137").Value = 729
$ws.Range("J137").Value = 3315.889
$ws.Range("K137").Value = 2187
$ws.Range("L137").Value = 9947.667000000001
$ws.Range("M137").Value = 363
$ws.Range("N137").Value = -15047.667
# Row 138
$ws.Range("H138").Value = 2590.9348
$ws.Range("I138").Value = 1981.6522
$ws.Range("K138").Value = 5944.9566
$ws.Range("M138").Value = -804.9565999999995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 7920
$ws.Range("I61").Value = 5316
$ws.Range("K61").Value = 5316
$ws.Range("M61").Value = -5104
# Row 74
$ws.Range("H74").Value = 22224786
$ws.Range("I74").Value = 25643676
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 25643676
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -25642802
$ws.Range("N74").Value = -3748
# Row 77
$ws.Range("H77").Value = 22224786
$ws.Range("I77").Value = 25643676
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 128218380
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -128214012
$ws.Range("N77").Value = -18736
# Row 122
$ws.Range("H122").Value = 2408.9473
$ws.Range("J122").Value = 3357.1428
$ws.Range("L122").Value = 10071.4284
$ws.Range("N122").Value = -14971.4284
# Row 136
$ws.Range("H136").Value = 7920
$ws.Range("I136").Value = 5316
$ws.Range("K136").Value = 15948
$ws.Range("M136").Value = -13398
# Row 139
$ws.Range("H139").Value = 68999.2
$ws.Range("J139").Value = 68999.2
$ws.Range("L139").Value = 68999.2
$ws.Range("N139").Value = -79279.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 108
$ws.Range("H108").Value = 96250
$ws.Range("J108").Value = 96250
$ws.Range("L108").Value = 96250
$ws.Range("N108").Value = -103930
# Row 123
$ws.Range("H123").Value = 77750
$ws.Range("J123").Value = 77750
$ws.Range("L123").Value = 77750
$ws.Range("N123").Value = -87550

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5951.4287
$ws.Range("I58").Value = 3024.5715
$ws.Range("K58").Value = 3024.5715
$ws.Range("M58").Value = -2821.5715
# Row 59
$ws.Range("H59").Value = 25000
$ws.Range("I59").Value = 25000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 25000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -23855
$ws.Range("N59").ClearContents()
# Row 122
$ws.Range("H122").Value = 4706.9033
$ws.Range("I122").Value = 1567.8125
$ws.Range("J122").Value = 8055.2666
$ws.Range("K122").Value = 4703.4375
$ws.Range("L122").Value = 24165.7998
$ws.Range("M122").Value = -2253.4375
$ws.Range("N122").Value = -29065.7998
# Row 123
$ws.Range("H123").Value = 36500
$ws.Range("J123").Value = 36500
$ws.Range("L123").Value = 36500
$ws.Range("N123").Value = -46300
# Row 130
$ws.Range("H130").Value = 37779
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
# Row 132
$ws.Range("H132").Value = 6049.1846
$ws.Range("I132").Value = 5888.4424
$ws.Range("J132").Value = 8500.5
$ws.Range("K132").Value = 17665.3272
$ws.Range("L132").Value = 25501.5
$ws.Range("M132").Value = -15135.3272
$ws.Range("N132").Value = -30561.5
# Row 134
$ws.Range("H134").Value = 5537.222
$ws.Range("I134").Value = 3853
$ws.Range("K134").Value = 11559
$ws.Range("M134").Value = -9024
# Row 136
$ws.Range("H136").Value = 5951.4287
$ws.Range("I136").Value = 3024.5715
$ws.Range("K136").Value = 9073.7145
$ws.Range("M136").Value = -6523.7145

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 2596.3572
$ws.Range("I122").Value = 494
$ws.Range("J122").Value = 2758.077
$ws.Range("K122").Value = 4446
$ws.Range("L122").Value = 24822.693
$ws.Range("M122").Value = -1996
$ws.Range("N122").Value = -29722.693

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 4109.067
$ws.Range("I126").Value = 3010.1667
$ws.Range("K126").Value = 9030.500100000001
$ws.Range("M126").Value = -6560.500100000001
# Row 132
$ws.Range("H132").Value = 57355.105
$ws.Range("I132").Value = 73979.92999999999
$ws.Range("K132").Value = 221939.79
$ws.Range("M132").Value = -219409.79

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 8087.25
$ws.Range("I40").Value = 7889.6875
$ws.Range("K40").Value = 7889.6875
$ws.Range("M40").Value = -7753.6875
# Row 68
$ws.Range("H68").Value = 4686.625
$ws.Range("I68").Value = 3023.85
$ws.Range("J68").Value = 13000.5
$ws.Range("K68").Value = 3023.85
$ws.Range("L68").Value = 13000.5
$ws.Range("M68").Value = -2274.85
$ws.Range("N68").Value = -14498.5
# Row 71
$ws.Range("H71").Value = 4686.625
$ws.Range("I71").Value = 3023.85
$ws.Range("J71").Value = 13000.5
$ws.Range("K71").Value = 15119.25
$ws.Range("L71").Value = 65002.5
$ws.Range("M71").Value = -11375.25
$ws.Range("N71").Value = -72490.5
# Row 122
$ws.Range("H122").Value = 8461.714
$ws.Range("J122").Value = 10302
$ws.Range("L122").Value = 30906
$ws.Range("N122").Value = -35806
# Row 132
$ws.Range("H132").Value = 7106.9
$ws.Range("I132").Value = 2333
$ws.Range("J132").Value = 9152.857
$ws.Range("K132").Value = 6999
$ws.Range("L132").Value = 27458.571
$ws.Range("M132").Value = -4469
$ws.Range("N132").Value = -32518.571

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1605.1428
$ws.Range("I100").Value = 556.75
$ws.Range("K100").Value = 1113.5
$ws.Range("M100").Value = -572.5
# Row 122
$ws.Range("H122").Value = 6085.278
$ws.Range("J122").Value = 11530.375
$ws.Range("L122").Value = 34591.125
$ws.Range("N122").Value = -39491.125
# Row 132
$ws.Range("H132").Value = 4436.3335
$ws.Range("I132").Value = 4443.5713
$ws.Range("J132").Value = 4411
$ws.Range("K132").Value = 13330.7139
$ws.Range("L132").Value = 13233
$ws.Range("M132").Value = -10800.7139
$ws.Range("N132").Value = -18293
